$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad / Changed date) for rows 2-6 from 45221 to 45224
$ws.Range("C2:C6").Value = 45224
